# "Add files via upload" -- the uploaded sheet shifted the phone-list rows:
# row 14 (+5511981446988) dropped out and every row from 15-34 moved up by
# one, with a brand new row landing at the bottom (row 34) and rows 19/20
# swapping which one carries the blank DDD cell.
# Only the cells whose value actually differs from the original are
# touched below (columns: A=Telefone, B=DDD, C=Data Inscricao), mirroring
# the upstream diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as literal text instead of
# re-parsing the leading "+" as a formula or the digits as a number (matches
# the source file, where every cell is a plain string).
function Set-Text($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

Set-Text "A14" "+553284691936"
Set-Text "B14" "32"

Set-Text "A15" "+5521978669363"
Set-Text "B15" "21"

Set-Text "A16" "+558781128650"
Set-Text "B16" "87"

Set-Text "A17" "+5511989613414"
Set-Text "B17" "11"
Set-Text "C17" "2024-10-31"

Set-Text "A18" "+5524981243416"
Set-Text "B18" "24"
Set-Text "C18" "2024-10-27"

Set-Text "A19" "+41766906567"
# DDD becomes blank for this row (still a text cell, not truly empty --
# mirrors the source file's self-closed <c t="inlineStr"/> pattern).
Set-Text "B19" ""

Set-Text "A20" "+553186055296"
Set-Text "B20" "31"
Set-Text "C20" "2024-10-26"

Set-Text "A21" "+556196229391"
Set-Text "B21" "61"

Set-Text "A22" "+559981971020"
Set-Text "B22" "99"
Set-Text "C22" "2024-10-24"

Set-Text "A23" "+5524999376084"
Set-Text "B23" "24"

Set-Text "A24" "+5511962253510"
Set-Text "B24" "11"

Set-Text "A25" "+5521966501459"
Set-Text "B25" "21"

Set-Text "A26" "+553384677059"
Set-Text "B26" "33"
Set-Text "C26" "2024-10-23"

Set-Text "A27" "+5521976781800"
Set-Text "B27" "21"

Set-Text "A28" "+5516993604875"
Set-Text "B28" "16"

Set-Text "A29" "+5521965496502"
Set-Text "B29" "21"
Set-Text "C29" "2024-10-22"

Set-Text "A30" "+5511972930265"
Set-Text "B30" "11"
Set-Text "C30" "2024-10-18"

Set-Text "A31" "+5511933453600"

Set-Text "A32" "+5511965520814"

Set-Text "A33" "+5521997899616"
Set-Text "B33" "21"

Set-Text "A34" "+5511950609656"
Set-Text "B34" "11"
Set-Text "C34" "2024-10-17"
